# W16 Thursday Commit 1
# Fills in the latest week of readings (7/28/20 - 7/31/20, serials 44037-44040)
# across all four tracking sheets, and adds a third leaf measurement column
# ("Leaf6"/"Leaf7"/"Leaf8") to the "Cardoon (2)" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Pansies Alive  (columns B:F = Pansy1..Pansy5)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("B28").Value = 5
$ws1.Range("C28").Value = 5
$ws1.Range("D28").Value = 4
$ws1.Range("E28").Value = 4
$ws1.Range("F28").Value = 0

$ws1.Range("B29").Value = 6
$ws1.Range("C29").Value = 8
$ws1.Range("D29").Value = 4
$ws1.Range("E29").Value = 4
$ws1.Range("F29").Value = 0

$ws1.Range("B30").Value = 6
$ws1.Range("C30").Value = 6
$ws1.Range("D30").Value = 4
$ws1.Range("E30").Value = 4
$ws1.Range("F30").Value = 0

$ws1.Range("A31").Value = 44040
$ws1.Range("A31").NumberFormat = "d-mmm-yy"
$ws1.Range("B31").Value = 5
$ws1.Range("C31").Value = 7
$ws1.Range("D31").Value = 4
$ws1.Range("E31").Value = 5
$ws1.Range("F31").Value = 0

# ---------------------------------------------------------------------
# Sheet 2: Pansies Dead  (columns B:F = Pansy1..Pansy5)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("B28").Value = 0
$ws2.Range("C28").Value = 1
$ws2.Range("D28").Value = 0
$ws2.Range("E28").Value = 0
$ws2.Range("F28").Value = 1

$ws2.Range("B29").Value = 0
$ws2.Range("C29").Value = 0
$ws2.Range("D29").Value = 0
$ws2.Range("E29").Value = 0
$ws2.Range("F29").Value = 1

$ws2.Range("B30").Value = 0
$ws2.Range("C30").Value = 1
$ws2.Range("D30").Value = 1
$ws2.Range("E30").Value = 0
$ws2.Range("F30").Value = 1

$ws2.Range("A31").Value = 44040
$ws2.Range("A31").NumberFormat = "d-mmm-yy"
$ws2.Range("B31").Value = 1
$ws2.Range("C31").Value = 0
$ws2.Range("D31").Value = 0
$ws2.Range("E31").Value = 0
$ws2.Range("F31").Value = 1

# ---------------------------------------------------------------------
# Sheet 3: Cardoon (1)  (columns C, F, G, H populated; B, D, E left blank
# for these rows, matching the existing pattern on this sheet)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("C28").Value = 19.75
$ws3.Range("F28").Value = 12
$ws3.Range("G28").Value = 17.5
$ws3.Range("H28").Value = 9

$ws3.Range("C29").Value = 20.25
$ws3.Range("F29").Value = 12
$ws3.Range("G29").Value = 17.5
$ws3.Range("H29").Value = 12

$ws3.Range("C30").Value = 19.5
$ws3.Range("F30").Value = 11
$ws3.Range("G30").Value = 17.5
$ws3.Range("H30").Value = 12

$ws3.Range("A31").Value = 44040
$ws3.Range("A31").NumberFormat = "d-mmm-yy"
$ws3.Range("C31").Value = 20.25
$ws3.Range("F31").Value = 11.5
$ws3.Range("G31").Value = 17.75
$ws3.Range("H31").Value = 12

# ---------------------------------------------------------------------
# Sheet 4: Cardoon (2)  -- gains a new column I ("Leaf8"); the old G/H
# headers are re-pointed at new shared-string entries "Leaf6"/"Leaf7"
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("G1").Value = "Leaf6"
$ws4.Range("H1").Value = "Leaf7"
$ws4.Range("I1").Value = "Leaf8"

$ws4.Range("D28").Value = 27
$ws4.Range("E28").Value = 24.5
$ws4.Range("F28").Value = 24.75
$ws4.Range("G28").Value = 24.75
$ws4.Range("H28").Value = 16.5

$ws4.Range("D29").Value = 27
$ws4.Range("E29").Value = 24
$ws4.Range("F29").Value = 24.5
$ws4.Range("G29").Value = 25.5
$ws4.Range("H29").Value = 17.45

$ws4.Range("D30").Value = 27
$ws4.Range("E30").Value = 24
$ws4.Range("F30").Value = 25
$ws4.Range("G30").Value = 25
$ws4.Range("H30").Value = 18.25
$ws4.Range("I30").Value = 1.25

$ws4.Range("D31").Value = 27
$ws4.Range("E31").Value = 24.5
$ws4.Range("F31").Value = 25.5
$ws4.Range("G31").Value = 25
$ws4.Range("H31").Value = 18.5
$ws4.Range("I31").Value = 2

# ---------------------------------------------------------------------
# View / selection state: final active sheet is Cardoon (2); leave a
# matching selection on each sheet before moving to the next one so the
# per-sheet <selection> element lines up with the authored state.
# ---------------------------------------------------------------------
$ws1.Range("I28").Select()
$ws2.Range("H30").Select()
$ws3.Range("M25").Select()
$ws4.Range("L27").Select()
$ws4.Activate()

$wb.Save()
